{"js": "// Fixed #418 Empty AQL expressions generate empty lines.\n// Remove the empty paragraph (no visible text) that sits right after the\n// \"Start of demonstration:\" paragraph.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Find the first paragraph whose text content is empty (the leftover blank\n// line produced by an empty AQL expression) and delete it entirely,\n// including its paragraph mark.\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.trim() === \"\") {\n    para.delete();\n    break;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Fixed #418 Empty AQL expressions generate empty lines.\n# Remove the empty paragraph (no visible text) that sits right after the\n# \"Start of demonstration:\" paragraph.\n\n$d = $word.ActiveDocument\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $para = $d.Paragraphs.Item($i)\n    $text = $para.Range.Text.Trim()\n    if ($text -eq \"\") {\n        $para.Range.Delete()\n        break\n    }\n}\n"}
